$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.169.68"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "1.905.93"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "'327.36"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "'0.4618"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("D8").Value = "'0.3940"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("D9").Value = "'46.71"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").Value = "'0.07910"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").Value = "'0.9963"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "'22.29"
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("D13").Value = "1.880.06"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'7.079"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").Value = "'5.739"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "'0.06960"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'88.56"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "'17.05"
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "29.192.49"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("D23").Value = "'5.328"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "2.171.38"
$ws.Range("E25").Value = "  +4.24%  "
$ws.Range("D26").Value = "'2.067"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "'156.72"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").Value = "'19.36"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "'6.052"
$ws.Range("E29").Value = "  +7.51%  "
$ws.Range("D30").Value = "'1.949"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "'118.08"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'0.09360"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "'0.9183"
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").Value = "'5.337"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").Value = "'1.351"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("D36").Value = "'3.283"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'1.198"
$ws.Range("E37").Value = "  +5.65%  "
$ws.Range("D38").Value = "'0.05816"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").Value = "'0.02099"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").Value = "'7.887"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'0.5725"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").Value = "'0.1790"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").Value = "'9.893"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "'2.302"
$ws.Range("E45").Value = "  +8.12%  "
$ws.Range("D46").Value = "'12.04"
$ws.Range("E46").Value = "  +3.98%  "
$ws.Range("D47").Value = "'0.5385"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("D48").Value = "'0.07043"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'1.863"
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("E50").Value = "  +5.59%  "
$ws.Range("D51").Value = "'113.01"
$ws.Range("E51").Value = "  +0.87%  "
